$wb = $excel.ActiveWorkbook

# ALC row 28
$ws = $wb.Worksheets.Item("ALC")
$F = $ws.Range("F28").Value2
$E = $ws.Range("E28").Value2
$newH = 497.33334
$newI = 570.86664
$newJ = 129.66667
$ws.Range("H28").Value2 = $newH
$ws.Range("I28").Value2 = $newI
$ws.Range("J28").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K28").Value2 = $K
$ws.Range("L28").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M28").Value2 = $E - $K
} else {
  $ws.Range("M28").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N28").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N28").ClearContents()
}

# ALC row 40
$ws = $wb.Worksheets.Item("ALC")
$F = $ws.Range("F40").Value2
$E = $ws.Range("E40").Value2
$newH = 14710796
$newI = 3476.2
$newJ = 20838846
$ws.Range("H40").Value2 = $newH
$ws.Range("I40").Value2 = $newI
$ws.Range("J40").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K40").Value2 = $K
$ws.Range("L40").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M40").Value2 = $E - $K
} else {
  $ws.Range("M40").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N40").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N40").ClearContents()
}

# ALC row 43
$ws = $wb.Worksheets.Item("ALC")
$F = $ws.Range("F43").Value2
$E = $ws.Range("E43").Value2
$newH = 3000
$newI = 1000
$newJ = 5000
$ws.Range("H43").Value2 = $newH
$ws.Range("I43").Value2 = $newI
$ws.Range("J43").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K43").Value2 = $K
$ws.Range("L43").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M43").Value2 = $E - $K
} else {
  $ws.Range("M43").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N43").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N43").ClearContents()
}

# ALC row 74
$ws = $wb.Worksheets.Item("ALC")
$F = $ws.Range("F74").Value2
$E = $ws.Range("E74").Value2
$newH = 12166.667
$newI = 6333.3335
$newJ = 14111.111
$ws.Range("H74").Value2 = $newH
$ws.Range("I74").Value2 = $newI
$ws.Range("J74").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K74").Value2 = $K
$ws.Range("L74").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M74").Value2 = $E - $K
} else {
  $ws.Range("M74").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N74").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N74").ClearContents()
}

# ALC row 77
$ws = $wb.Worksheets.Item("ALC")
$F = $ws.Range("F77").Value2
$E = $ws.Range("E77").Value2
$newH = 12166.667
$newI = 6333.3335
$newJ = 14111.111
$ws.Range("H77").Value2 = $newH
$ws.Range("I77").Value2 = $newI
$ws.Range("J77").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K77").Value2 = $K
$ws.Range("L77").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M77").Value2 = $E - $K
} else {
  $ws.Range("M77").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N77").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N77").ClearContents()
}

# ALC row 80
$ws = $wb.Worksheets.Item("ALC")
$F = $ws.Range("F80").Value2
$E = $ws.Range("E80").Value2
$newH = 8366.799999999999
$newI = 6869.933
$newJ = 10612.1
$ws.Range("H80").Value2 = $newH
$ws.Range("I80").Value2 = $newI
$ws.Range("J80").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K80").Value2 = $K
$ws.Range("L80").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M80").Value2 = $E - $K
} else {
  $ws.Range("M80").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N80").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N80").ClearContents()
}

# ALC row 83
$ws = $wb.Worksheets.Item("ALC")
$F = $ws.Range("F83").Value2
$E = $ws.Range("E83").Value2
$newH = 8366.799999999999
$newI = 6869.933
$newJ = 10612.1
$ws.Range("H83").Value2 = $newH
$ws.Range("I83").Value2 = $newI
$ws.Range("J83").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K83").Value2 = $K
$ws.Range("L83").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M83").Value2 = $E - $K
} else {
  $ws.Range("M83").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N83").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N83").ClearContents()
}

# ALC row 86
$ws = $wb.Worksheets.Item("ALC")
$F = $ws.Range("F86").Value2
$E = $ws.Range("E86").Value2
$newH = 1814.1428
$newI = 1866.6666
$newJ = 1499
$ws.Range("H86").Value2 = $newH
$ws.Range("I86").Value2 = $newI
$ws.Range("J86").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K86").Value2 = $K
$ws.Range("L86").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M86").Value2 = $E - $K
} else {
  $ws.Range("M86").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N86").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N86").ClearContents()
}

# ALC row 89
$ws = $wb.Worksheets.Item("ALC")
$F = $ws.Range("F89").Value2
$E = $ws.Range("E89").Value2
$newH = 1814.1428
$newI = 1866.6666
$newJ = 1499
$ws.Range("H89").Value2 = $newH
$ws.Range("I89").Value2 = $newI
$ws.Range("J89").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K89").Value2 = $K
$ws.Range("L89").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M89").Value2 = $E - $K
} else {
  $ws.Range("M89").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N89").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N89").ClearContents()
}

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$F = $ws.Range("F112").Value2
$E = $ws.Range("E112").Value2
$newH = 1830
$newI = 1830
$newJ = 0
$ws.Range("H112").Value2 = $newH
$ws.Range("I112").Value2 = $newI
$ws.Range("J112").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K112").Value2 = $K
$ws.Range("L112").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M112").Value2 = $E - $K
} else {
  $ws.Range("M112").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N112").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N112").ClearContents()
}

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$F = $ws.Range("F137").Value2
$E = $ws.Range("E137").Value2
$newH = 7250.591
$newI = 7349.1333
$newJ = 7039.4287
$ws.Range("H137").Value2 = $newH
$ws.Range("I137").Value2 = $newI
$ws.Range("J137").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K137").Value2 = $K
$ws.Range("L137").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M137").Value2 = $E - $K
} else {
  $ws.Range("M137").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N137").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N137").ClearContents()
}

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$F = $ws.Range("F138").Value2
$E = $ws.Range("E138").Value2
$newH = 7131.515
$newI = 6942.6
$newJ = 7422.154
$ws.Range("H138").Value2 = $newH
$ws.Range("I138").Value2 = $newI
$ws.Range("J138").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K138").Value2 = $K
$ws.Range("L138").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M138").Value2 = $E - $K
} else {
  $ws.Range("M138").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N138").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N138").ClearContents()
}

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$F = $ws.Range("F32").Value2
$E = $ws.Range("E32").Value2
$newH = 10523.2
$newI = 9024.713
$newJ = 34499
$ws.Range("H32").Value2 = $newH
$ws.Range("I32").Value2 = $newI
$ws.Range("J32").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K32").Value2 = $K
$ws.Range("L32").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M32").Value2 = $E - $K
} else {
  $ws.Range("M32").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N32").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N32").ClearContents()
}

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$F = $ws.Range("F45").Value2
$E = $ws.Range("E45").Value2
$newH = 83109.03999999999
$newI = 113201.555
$newJ = 5728.2856
$ws.Range("H45").Value2 = $newH
$ws.Range("I45").Value2 = $newI
$ws.Range("J45").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K45").Value2 = $K
$ws.Range("L45").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M45").Value2 = $E - $K
} else {
  $ws.Range("M45").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N45").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N45").ClearContents()
}

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$F = $ws.Range("F63").Value2
$E = $ws.Range("E63").Value2
$newH = 3733.3333
$newI = 3733.3333
$newJ = 0
$ws.Range("H63").Value2 = $newH
$ws.Range("I63").Value2 = $newI
$ws.Range("J63").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K63").Value2 = $K
$ws.Range("L63").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M63").Value2 = $E - $K
} else {
  $ws.Range("M63").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N63").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N63").ClearContents()
}

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$F = $ws.Range("F66").Value2
$E = $ws.Range("E66").Value2
$newH = 3733.3333
$newI = 3733.3333
$newJ = 0
$ws.Range("H66").Value2 = $newH
$ws.Range("I66").Value2 = $newI
$ws.Range("J66").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K66").Value2 = $K
$ws.Range("L66").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M66").Value2 = $E - $K
} else {
  $ws.Range("M66").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N66").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N66").ClearContents()
}

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$F = $ws.Range("F74").Value2
$E = $ws.Range("E74").Value2
$newH = 811.2222
$newI = 811.2222
$newJ = 0
$ws.Range("H74").Value2 = $newH
$ws.Range("I74").Value2 = $newI
$ws.Range("J74").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K74").Value2 = $K
$ws.Range("L74").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M74").Value2 = $E - $K
} else {
  $ws.Range("M74").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N74").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N74").ClearContents()
}

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$F = $ws.Range("F77").Value2
$E = $ws.Range("E77").Value2
$newH = 811.2222
$newI = 811.2222
$newJ = 0
$ws.Range("H77").Value2 = $newH
$ws.Range("I77").Value2 = $newI
$ws.Range("J77").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K77").Value2 = $K
$ws.Range("L77").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M77").Value2 = $E - $K
} else {
  $ws.Range("M77").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N77").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N77").ClearContents()
}

# ARM row 80
$ws = $wb.Worksheets.Item("ARM")
$F = $ws.Range("F80").Value2
$E = $ws.Range("E80").Value2
$newH = 40936
$newI = 0
$newJ = 40936
$ws.Range("H80").Value2 = $newH
$ws.Range("I80").Value2 = $newI
$ws.Range("J80").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K80").Value2 = $K
$ws.Range("L80").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M80").Value2 = $E - $K
} else {
  $ws.Range("M80").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N80").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N80").ClearContents()
}

# ARM row 83
$ws = $wb.Worksheets.Item("ARM")
$F = $ws.Range("F83").Value2
$E = $ws.Range("E83").Value2
$newH = 40936
$newI = 0
$newJ = 40936
$ws.Range("H83").Value2 = $newH
$ws.Range("I83").Value2 = $newI
$ws.Range("J83").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K83").Value2 = $K
$ws.Range("L83").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M83").Value2 = $E - $K
} else {
  $ws.Range("M83").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N83").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N83").ClearContents()
}

# ARM row 88
$ws = $wb.Worksheets.Item("ARM")
$F = $ws.Range("F88").Value2
$E = $ws.Range("E88").Value2
$newH = 3696.5833
$newI = 3885
$newJ = 3562
$ws.Range("H88").Value2 = $newH
$ws.Range("I88").Value2 = $newI
$ws.Range("J88").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K88").Value2 = $K
$ws.Range("L88").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M88").Value2 = $E - $K
} else {
  $ws.Range("M88").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N88").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N88").ClearContents()
}

# ARM row 91
$ws = $wb.Worksheets.Item("ARM")
$F = $ws.Range("F91").Value2
$E = $ws.Range("E91").Value2
$newH = 3696.5833
$newI = 3885
$newJ = 3562
$ws.Range("H91").Value2 = $newH
$ws.Range("I91").Value2 = $newI
$ws.Range("J91").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K91").Value2 = $K
$ws.Range("L91").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M91").Value2 = $E - $K
} else {
  $ws.Range("M91").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N91").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N91").ClearContents()
}

# ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$F = $ws.Range("F97").Value2
$E = $ws.Range("E97").Value2
$newH = 911.76666
$newI = 974.5599999999999
$newJ = 597.8
$ws.Range("H97").Value2 = $newH
$ws.Range("I97").Value2 = $newI
$ws.Range("J97").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K97").Value2 = $K
$ws.Range("L97").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M97").Value2 = $E - $K
} else {
  $ws.Range("M97").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N97").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N97").ClearContents()
}

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$F = $ws.Range("F122").Value2
$E = $ws.Range("E122").Value2
$newH = 2072.2942
$newI = 2024.3077
$newJ = 2228.25
$ws.Range("H122").Value2 = $newH
$ws.Range("I122").Value2 = $newI
$ws.Range("J122").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K122").Value2 = $K
$ws.Range("L122").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M122").Value2 = $E - $K
} else {
  $ws.Range("M122").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N122").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N122").ClearContents()
}

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$F = $ws.Range("F132").Value2
$E = $ws.Range("E132").Value2
$newH = 2717.8438
$newI = 2782.5356
$newJ = 2265
$ws.Range("H132").Value2 = $newH
$ws.Range("I132").Value2 = $newI
$ws.Range("J132").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K132").Value2 = $K
$ws.Range("L132").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M132").Value2 = $E - $K
} else {
  $ws.Range("M132").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N132").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N132").ClearContents()
}

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$F = $ws.Range("F94").Value2
$E = $ws.Range("E94").Value2
$newH = 265.33334
$newI = 265.33334
$newJ = 0
$ws.Range("H94").Value2 = $newH
$ws.Range("I94").Value2 = $newI
$ws.Range("J94").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K94").Value2 = $K
$ws.Range("L94").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M94").Value2 = $E - $K
} else {
  $ws.Range("M94").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N94").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N94").ClearContents()
}

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$F = $ws.Range("F99").Value2
$E = $ws.Range("E99").Value2
$newH = 1408.6154
$newI = 931.4
$newJ = 2999.3333
$ws.Range("H99").Value2 = $newH
$ws.Range("I99").Value2 = $newI
$ws.Range("J99").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K99").Value2 = $K
$ws.Range("L99").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M99").Value2 = $E - $K
} else {
  $ws.Range("M99").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N99").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N99").ClearContents()
}

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$F = $ws.Range("F134").Value2
$E = $ws.Range("E134").Value2
$newH = 3575.1924
$newI = 3040.913
$newJ = 7671.3335
$ws.Range("H134").Value2 = $newH
$ws.Range("I134").Value2 = $newI
$ws.Range("J134").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K134").Value2 = $K
$ws.Range("L134").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M134").Value2 = $E - $K
} else {
  $ws.Range("M134").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N134").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N134").ClearContents()
}

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$F = $ws.Range("F31").Value2
$E = $ws.Range("E31").Value2
$newH = 5585.6
$newI = 3588.5
$newJ = 9579.799999999999
$ws.Range("H31").Value2 = $newH
$ws.Range("I31").Value2 = $newI
$ws.Range("J31").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K31").Value2 = $K
$ws.Range("L31").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M31").Value2 = $E - $K
} else {
  $ws.Range("M31").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N31").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N31").ClearContents()
}

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$F = $ws.Range("F34").Value2
$E = $ws.Range("E34").Value2
$newH = 5585.6
$newI = 3588.5
$newJ = 9579.799999999999
$ws.Range("H34").Value2 = $newH
$ws.Range("I34").Value2 = $newI
$ws.Range("J34").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K34").Value2 = $K
$ws.Range("L34").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M34").Value2 = $E - $K
} else {
  $ws.Range("M34").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N34").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N34").ClearContents()
}

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$F = $ws.Range("F99").Value2
$E = $ws.Range("E99").Value2
$newH = 7957.154
$newI = 6112
$newJ = 8777.223
$ws.Range("H99").Value2 = $newH
$ws.Range("I99").Value2 = $newI
$ws.Range("J99").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K99").Value2 = $K
$ws.Range("L99").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M99").Value2 = $E - $K
} else {
  $ws.Range("M99").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N99").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N99").ClearContents()
}

# CRP row 105
$ws = $wb.Worksheets.Item("CRP")
$F = $ws.Range("F105").Value2
$E = $ws.Range("E105").Value2
$newH = 1594.7778
$newI = 1282.44
$newJ = 5499
$ws.Range("H105").Value2 = $newH
$ws.Range("I105").Value2 = $newI
$ws.Range("J105").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K105").Value2 = $K
$ws.Range("L105").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M105").Value2 = $E - $K
} else {
  $ws.Range("M105").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N105").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N105").ClearContents()
}

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$F = $ws.Range("F126").Value2
$E = $ws.Range("E126").Value2
$newH = 7957.154
$newI = 6112
$newJ = 8777.223
$ws.Range("H126").Value2 = $newH
$ws.Range("I126").Value2 = $newI
$ws.Range("J126").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K126").Value2 = $K
$ws.Range("L126").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M126").Value2 = $E - $K
} else {
  $ws.Range("M126").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N126").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N126").ClearContents()
}

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$F = $ws.Range("F132").Value2
$E = $ws.Range("E132").Value2
$newH = 373623.44
$newI = 557424.1
$newJ = 6022.1113
$ws.Range("H132").Value2 = $newH
$ws.Range("I132").Value2 = $newI
$ws.Range("J132").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K132").Value2 = $K
$ws.Range("L132").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M132").Value2 = $E - $K
} else {
  $ws.Range("M132").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N132").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N132").ClearContents()
}

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$F = $ws.Range("F134").Value2
$E = $ws.Range("E134").Value2
$newH = 3589.3438
$newI = 1979.3636
$newJ = 7131.3
$ws.Range("H134").Value2 = $newH
$ws.Range("I134").Value2 = $newI
$ws.Range("J134").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K134").Value2 = $K
$ws.Range("L134").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M134").Value2 = $E - $K
} else {
  $ws.Range("M134").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N134").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N134").ClearContents()
}

# CRP row 141
$ws = $wb.Worksheets.Item("CRP")
$F = $ws.Range("F141").Value2
$E = $ws.Range("E141").Value2
$newH = 239362.66
$newI = 50000
$newJ = 272779.6
$ws.Range("H141").Value2 = $newH
$ws.Range("I141").Value2 = $newI
$ws.Range("J141").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K141").Value2 = $K
$ws.Range("L141").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M141").Value2 = $E - $K
} else {
  $ws.Range("M141").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N141").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N141").ClearContents()
}

# CUL row 2
$ws = $wb.Worksheets.Item("CUL")
$F = $ws.Range("F2").Value2
$E = $ws.Range("E2").Value2
$newH = 1536.2084
$newI = 86.42856999999999
$newJ = 11684.667
$ws.Range("H2").Value2 = $newH
$ws.Range("I2").Value2 = $newI
$ws.Range("J2").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K2").Value2 = $K
$ws.Range("L2").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M2").Value2 = $E - $K
} else {
  $ws.Range("M2").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N2").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N2").ClearContents()
}

# CUL row 37
$ws = $wb.Worksheets.Item("CUL")
$F = $ws.Range("F37").Value2
$E = $ws.Range("E37").Value2
$newH = 164804.44
$newI = 0
$newJ = 164804.44
$ws.Range("H37").Value2 = $newH
$ws.Range("I37").Value2 = $newI
$ws.Range("J37").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K37").Value2 = $K
$ws.Range("L37").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M37").Value2 = $E - $K
} else {
  $ws.Range("M37").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N37").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N37").ClearContents()
}

# CUL row 38
$ws = $wb.Worksheets.Item("CUL")
$F = $ws.Range("F38").Value2
$E = $ws.Range("E38").Value2
$newH = 464.4
$newI = 20
$newJ = 575.5
$ws.Range("H38").Value2 = $newH
$ws.Range("I38").Value2 = $newI
$ws.Range("J38").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K38").Value2 = $K
$ws.Range("L38").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M38").Value2 = $E - $K
} else {
  $ws.Range("M38").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N38").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N38").ClearContents()
}

# CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$F = $ws.Range("F122").Value2
$E = $ws.Range("E122").Value2
$newH = 4857.794
$newI = 1773.25
$newJ = 5806.885
$ws.Range("H122").Value2 = $newH
$ws.Range("I122").Value2 = $newI
$ws.Range("J122").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K122").Value2 = $K
$ws.Range("L122").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M122").Value2 = $E - $K
} else {
  $ws.Range("M122").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N122").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N122").ClearContents()
}

# CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$F = $ws.Range("F132").Value2
$E = $ws.Range("E132").Value2
$newH = 100002400
$newI = 0
$newJ = 100002400
$ws.Range("H132").Value2 = $newH
$ws.Range("I132").Value2 = $newI
$ws.Range("J132").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K132").Value2 = $K
$ws.Range("L132").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M132").Value2 = $E - $K
} else {
  $ws.Range("M132").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N132").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N132").ClearContents()
}

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$F = $ws.Range("F122").Value2
$E = $ws.Range("E122").Value2
$newH = 2676.087
$newI = 2761.2222
$newJ = 2369.6
$ws.Range("H122").Value2 = $newH
$ws.Range("I122").Value2 = $newI
$ws.Range("J122").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K122").Value2 = $K
$ws.Range("L122").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M122").Value2 = $E - $K
} else {
  $ws.Range("M122").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N122").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N122").ClearContents()
}

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$F = $ws.Range("F132").Value2
$E = $ws.Range("E132").Value2
$newH = 1882.1034
$newI = 1967.1538
$newJ = 1145
$ws.Range("H132").Value2 = $newH
$ws.Range("I132").Value2 = $newI
$ws.Range("J132").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K132").Value2 = $K
$ws.Range("L132").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M132").Value2 = $E - $K
} else {
  $ws.Range("M132").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N132").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N132").ClearContents()
}

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$F = $ws.Range("F46").Value2
$E = $ws.Range("E46").Value2
$newH = 1699.75
$newI = 1666.3334
$newJ = 1800
$ws.Range("H46").Value2 = $newH
$ws.Range("I46").Value2 = $newI
$ws.Range("J46").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K46").Value2 = $K
$ws.Range("L46").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M46").Value2 = $E - $K
} else {
  $ws.Range("M46").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N46").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N46").ClearContents()
}

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$F = $ws.Range("F132").Value2
$E = $ws.Range("E132").Value2
$newH = 208881.98
$newI = 250642.69
$newJ = 13998.667
$ws.Range("H132").Value2 = $newH
$ws.Range("I132").Value2 = $newI
$ws.Range("J132").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K132").Value2 = $K
$ws.Range("L132").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M132").Value2 = $E - $K
} else {
  $ws.Range("M132").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N132").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N132").ClearContents()
}

# WVR row 46
$ws = $wb.Worksheets.Item("WVR")
$F = $ws.Range("F46").Value2
$E = $ws.Range("E46").Value2
$newH = 79992.5
$newI = 0
$newJ = 79992.5
$ws.Range("H46").Value2 = $newH
$ws.Range("I46").Value2 = $newI
$ws.Range("J46").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K46").Value2 = $K
$ws.Range("L46").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M46").Value2 = $E - $K
} else {
  $ws.Range("M46").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N46").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N46").ClearContents()
}

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$F = $ws.Range("F132").Value2
$E = $ws.Range("E132").Value2
$newH = 136302.62
$newI = 178134.28
$newJ = 3835.7222
$ws.Range("H132").Value2 = $newH
$ws.Range("I132").Value2 = $newI
$ws.Range("J132").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K132").Value2 = $K
$ws.Range("L132").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M132").Value2 = $E - $K
} else {
  $ws.Range("M132").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N132").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N132").ClearContents()
}

# WVR row 134
$ws = $wb.Worksheets.Item("WVR")
$F = $ws.Range("F134").Value2
$E = $ws.Range("E134").Value2
$newH = 79992.5
$newI = 0
$newJ = 79992.5
$ws.Range("H134").Value2 = $newH
$ws.Range("I134").Value2 = $newI
$ws.Range("J134").Value2 = $newJ
$K = $newI * $F
$L = $newJ * $F
$ws.Range("K134").Value2 = $K
$ws.Range("L134").Value2 = $L
if ($newI -ne 0) {
  $ws.Range("M134").Value2 = $E - $K
} else {
  $ws.Range("M134").ClearContents()
}
if ($newJ -ne 0) {
  $ws.Range("N134").Value2 = (-2 * $E) - $L
} else {
  $ws.Range("N134").ClearContents()
}
